# "Add cantrals by cantons" - restructure AI 2013 sheet:
#  - Collapse the old two-row header (row1 "units row" + row2 "labels row")
#    into a single header row with explicit column titles.
#  - Shift the two data rows up by one.
#  - Drop the now-superfluous trailing blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the two existing data rows before we start overwriting cells ---
# (.Value2 is used for reads -- .Value round-trips as an opaque COM stub here)
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2
$f3 = $ws.Range("F3").Value2
$g3 = $ws.Range("G3").Value2
$h3 = $ws.Range("H3").Value2
$i3 = $ws.Range("I3").Value2
$j3 = $ws.Range("J3").Value2
$k3 = $ws.Range("K3").Value2

$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$c4 = $ws.Range("C4").Value2
$d4 = $ws.Range("D4").Value2
$e4 = $ws.Range("E4").Value2
$f4 = $ws.Range("F4").Value2
$g4 = $ws.Range("G4").Value2
$h4 = $ws.Range("H4").Value2
$i4 = $ws.Range("I4").Value2
$j4 = $ws.Range("J4").Value2
$k4 = $ws.Range("K4").Value2

# --- a transient named style == Arial 9, General format, font-only applied ---
# (matches the pre-existing "label" font used elsewhere in the sheet, fontId 2)
$wb.Styles.Add("__hdr")
$hdrStyle = $wb.Styles.Item("__hdr")
$hdrStyle.Font.Name = "Arial"
$hdrStyle.Font.Size = 9

# --- row 1: single header row ---
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"
$ws.Range("F1:K1").Style = "__hdr"

# the helper named style has done its job -- drop it so it doesn't linger
# as an extra entry; the cellXf it produced stays attached to F1:K1.
$wb.Styles.Item("__hdr").Delete()

# --- row 2: first data row (was row 3) ---
$ws.Range("A2").Value = $a3
$ws.Range("B2").Value = $b3
$ws.Range("C2").Value = $c3
$ws.Range("D2").Value = $d3
$ws.Range("E2").Value = $e3
$ws.Range("F2").Value = $f3
$ws.Range("G2").Value = $g3
$ws.Range("H2").Value = $h3
$ws.Range("I2").Value = $i3
$ws.Range("J2").Value = $j3
$ws.Range("K2").Value = $k3

# --- row 3: second data row (was row 4) ---
$ws.Range("A3").Value = $a4
$ws.Range("B3").Value = $b4
$ws.Range("C3").Value = $c4
$ws.Range("D3").Value = $d4
$ws.Range("E3").Value = $e4
$ws.Range("F3").Value = $f4
$ws.Range("G3").Value = $g4
$ws.Range("H3").Value = $h4
$ws.Range("I3").Value = $i4
$ws.Range("J3").Value = $j4
$ws.Range("K3").Value = $k4

# --- row 4 reverts to an empty placeholder row (matches the old row 5 pattern) ---
$ws.Range("A4").Value = ""

# the sheet had one blank trailing row too many -- remove it
$ws.Rows(78).Delete()

# mirror the author's final selection
$ws.Range("A2:K2").Select()
